# Add a new "FRA" column (H) to the JudgeDashboard sheet with Y/N/U flags
# per row, matching the values added to the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("H1").Value = "FRA"

# Data rows 2-23 (values correspond to shared strings Y / N / U)
$values = @{
    2  = "Y"
    3  = "Y"
    4  = "Y"
    5  = "N"
    6  = "N"
    7  = "Y"
    8  = "Y"
    9  = "Y"
    10 = "U"
    11 = "N"
    12 = "N"
    13 = "U"
    14 = "Y"
    15 = "Y"
    16 = "Y"
    17 = "Y"
    18 = "Y"
    19 = "Y"
    20 = "Y"
    21 = "N"
    22 = "N"
    23 = "N"
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Cells.Item($row, 8).Value = $values[$row]
}

# Touch the trailing spacer row so the sheet's used range/dimension and row
# span extend to include the new column, matching the widened data grid.
$ws.Cells.Item(24, 8).Font.Bold = $false

# Update the selection/active cell to the last populated cell in the new column
$ws.Range("H24").Select() | Out-Null

Write-Output "done"
